$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.533508
$ws.Range("H2").Value = 10.600524
$ws.Range("I2").Value = 0.0236509744414791
$ws.Range("J2").Value = 0.02365097444147911
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2113696666666667
$ws.Range("N2").Value = 0.634109
$ws.Range("O2").Value = 0.03795977003925348
$ws.Range("P2").Value = 0.03795977003925347
$ws.Range("Q2").Value = 0.746876408124
$ws.Range("R2").Value = 6.721887673116
$ws.Range("S2").Value = 0.0008977855510028083
$ws.Range("T2").Value = 0.0008977855510028083

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.533508
$ws.Range("H3").Value = 10.600524
$ws.Range("I3").Value = 0.0236509744414791
$ws.Range("J3").Value = 0.02365097444147911
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.589504333333333
$ws.Range("N3").Value = 4.768513
$ws.Range("O3").Value = 0.2854582680725092
$ws.Range("P3").Value = 0.2854582680725091
$ws.Range("Q3").Value = 5.616526277868
$ws.Range("R3").Value = 50.54873650081201
$ws.Range("S3").Value = 0.006751366202291806
$ws.Range("T3").Value = 0.006751366202291805

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.533508
$ws.Range("H4").Value = 10.600524
$ws.Range("I4").Value = 0.0236509744414791
$ws.Range("J4").Value = 0.02365097444147911
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.767380666666666
$ws.Range("N4").Value = 11.302142
$ws.Range("O4").Value = 0.6765819618882374
$ws.Range("P4").Value = 0.6765819618882374
$ws.Range("Q4").Value = 13.312069724712
$ws.Range("R4").Value = 119.808627522408
$ws.Range("S4").Value = 0.01600182268818449
$ws.Range("T4").Value = 0.01600182268818449

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 133.9582823333334
$ws.Range("H5").Value = 401.874847
$ws.Range("I5").Value = 0.8966284812968046
$ws.Range("J5").Value = 0.8966284812968046
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2113696666666667
$ws.Range("N5").Value = 0.634109
$ws.Range("O5").Value = 0.03795977003925348
$ws.Range("P5").Value = 0.03795977003925347
$ws.Range("Q5").Value = 28.3147174840359
$ws.Range("R5").Value = 254.832457356323
$ws.Range("S5").Value = 0.03403581096067179
$ws.Range("T5").Value = 0.03403581096067178

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 133.9582823333334
$ws.Range("H6").Value = 401.874847
$ws.Range("I6").Value = 0.8966284812968046
$ws.Range("J6").Value = 0.8966284812968046
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.589504333333333
$ws.Range("N6").Value = 4.768513
$ws.Range("O6").Value = 0.2854582680725092
$ws.Range("P6").Value = 0.2854582680725091
$ws.Range("Q6").Value = 212.9272702547235
$ws.Range("R6").Value = 1916.345432292511
$ws.Range("S6").Value = 0.25595001337547
$ws.Range("T6").Value = 0.25595001337547

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 133.9582823333334
$ws.Range("H7").Value = 401.874847
$ws.Range("I7").Value = 0.8966284812968046
$ws.Range("J7").Value = 0.8966284812968046
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.767380666666666
$ws.Range("N7").Value = 11.302142
$ws.Range("O7").Value = 0.6765819618882374
$ws.Range("P7").Value = 0.6765819618882374
$ws.Range("Q7").Value = 504.6718430024749
$ws.Range("R7").Value = 4542.046587022274
$ws.Range("S7").Value = 0.6066426569606628
$ws.Range("T7").Value = 0.6066426569606628

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.910426
$ws.Range("H8").Value = 35.731278
$ws.Range("I8").Value = 0.07972054426171619
$ws.Range("J8").Value = 0.07972054426171619
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2113696666666667
$ws.Range("N8").Value = 0.634109
$ws.Range("O8").Value = 0.03795977003925348
$ws.Range("P8").Value = 0.03795977003925347
$ws.Range("Q8").Value = 2.517502773478
$ws.Range("R8").Value = 22.657524961302
$ws.Range("S8").Value = 0.003026173527578875
$ws.Range("T8").Value = 0.003026173527578875

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.910426
$ws.Range("H9").Value = 35.731278
$ws.Range("I9").Value = 0.07972054426171619
$ws.Range("J9").Value = 0.07972054426171619
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.589504333333333
$ws.Range("N9").Value = 4.768513
$ws.Range("O9").Value = 0.2854582680725092
$ws.Range("P9").Value = 0.2854582680725091
$ws.Range("Q9").Value = 18.931673738846
$ws.Range("R9").Value = 170.385063649614
$ws.Range("S9").Value = 0.02275688849474732
$ws.Range("T9").Value = 0.02275688849474731

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.910426
$ws.Range("H10").Value = 35.731278
$ws.Range("I10").Value = 0.07972054426171619
$ws.Range("J10").Value = 0.07972054426171619
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.767380666666666
$ws.Range("N10").Value = 11.302142
$ws.Range("O10").Value = 0.6765819618882374
$ws.Range("P10").Value = 0.6765819618882374
$ws.Range("Q10").Value = 44.871108644164
$ws.Range("R10").Value = 403.839977797476
$ws.Range("S10").Value = 0.05393748223939001
$ws.Range("T10").Value = 0.05393748223939001

Write-Output "done"